$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The [Config] tag was renamed to [ExcelLENT]
$ws.Range("A2").Value = "[ExcelLENT]"

# Move the active selection to A2 (was A8)
$ws.Range("A2").Select()
